$wb = $excel.ActiveWorkbook

# --- survey sheet: add a new select_one_with_other question (colors) ---
$survey = $wb.Worksheets.Item("survey")

# Row 14: select_one_with_other colors / color / What is your favorite color?
$survey.Range("B14").Value = "select_one_with_other colors"
$survey.Range("E14").Value = "color"
$survey.Range("F14").Value = "What is your favorite color?"

# Row 15: a note warning about `selected()` with an argument not in the choice list
$survey.Range("A15").Value = "selected function with arguement not included in choices."
$survey.Range("B15").Value = "note"
$survey.Range("D15").Value = "selected(data('color'), 'teal')"
$survey.Range("F15").Value = "Teal is a good choice."

# --- choices sheet: add the "colors" choice list (red/green/blue) ---
$choices = $wb.Worksheets.Item("choices")

$choices.Range("A15").Value = "colors"
$choices.Range("B15").Value = "red"
$choices.Range("D15").Value = "Red"

$choices.Range("A16").Value = "colors"
$choices.Range("B16").Value = "green"
$choices.Range("D16").Value = "Green"

$choices.Range("A17").Value = "colors"
$choices.Range("B17").Value = "blue"
$choices.Range("D17").Value = "Blue"
